$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) for rows 2 through 17
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25)
for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
